$d = $word.ActiveDocument

# The document ends with three empty "Normal" paragraphs (right after the
# "easy to use Ios app for electronic raffle tickets" project entry).
# We insert the new "INTERNSHIP PROJECTS" section right after the first
# of those trailing empty paragraphs, keeping the last two trailing
# empty paragraphs intact at the end of the document.

$anchorIndex = $d.Paragraphs.Count - 2

$lines = @(
    "",
    "INTERNSHIP PROJECTS",
    "",
    "Roadsign classifier",
    "",
    "ported python tensorflow graph to c++",
    "",
    "optimised tensorflow graph using tensorrt",
    "",
    "created a data augmentation script",
    "",
    "created a data sanitisation tool to reduce false positives",
    "",
    "created various scripts to clean up the BDD100k dataset"
)

$insertAfterIndex = $anchorIndex
foreach ($line in $lines) {
    $ref = $d.Paragraphs.Item($insertAfterIndex).Range
    $ref.InsertParagraphAfter()
    $insertAfterIndex = $insertAfterIndex + 1
    if ($line -ne "") {
        $newPara = $d.Paragraphs.Item($insertAfterIndex)
        $newPara.Range.InsertAfter($line)
    }
}
